# "Logged Week 15 and simulated Week 16"
# Appends one new week's worth of per-play/per-game numbers to the
# running logs on YDS and ST, and updates the season-to-date totals on
# OFF, DEF, ST, TURNS and PEN.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: append this week's run/pass yardage-by-play logs
# ---------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value2 = $ydsWs.Range("B2").Value2 + " 4 15 3 15 4 3 6 2 4 1 2 7 0 9 4 6 8 4 3 1 8 2 4 8 22 1 1 3 0 2 2 3 2 6 15 2 2"
$ydsWs.Range("C2").Value2 = $ydsWs.Range("C2").Value2 + " 5 11 9 7 0 5 0 -1 3 3 1 2 3 0 1 5 0 1 32"
$ydsWs.Range("B3").Value2 = $ydsWs.Range("B3").Value2 + " 9 4 12 12 13 15 4 2 27 7 5 4 17 6 29 8 5 15 15 8 5 14"
$ydsWs.Range("C3").Value2 = $ydsWs.Range("C3").Value2 + " 12 14 5 9 20 11 12 -1 17 4 3 6 5 23 10 40 5 3 9 0 6 69 1 5 10 20 7 10 27 14 34"

# ---------------------------------------------------------------
# OFF sheet: season totals through this week
# ---------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 381
$offWs.Range("D2").Value = 26
$offWs.Range("F2").Value = 99
$offWs.Range("G2").Value = 110
$offWs.Range("J2").Value = 50
$offWs.Range("L2").Value = 647
$offWs.Range("M2").Value = 428
$offWs.Range("O2").Value = 50
$offWs.Range("P2").Value = 28
$offWs.Range("Q2").Value = 1101

$offWs.Range("C3").Value = 351
$offWs.Range("D3").Value = 10
$offWs.Range("E3").Value = 71
$offWs.Range("F3").Value = 247
$offWs.Range("G3").Value = 64
$offWs.Range("H3").Value = 59
$offWs.Range("I3").Value = 115
$offWs.Range("J3").Value = 119

# ---------------------------------------------------------------
# DEF sheet: season totals through this week
# ---------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value = 6
$defWs.Range("C2").Value = 384
$defWs.Range("D2").Value = 26
$defWs.Range("F2").Value = 114
$defWs.Range("G2").Value = 122
$defWs.Range("J2").Value = 65
$defWs.Range("L2").Value = 547
$defWs.Range("M2").Value = 335
$defWs.Range("O2").Value = 32
$defWs.Range("Q2").Value = 1003

$defWs.Range("C3").Value = 310
$defWs.Range("D3").Value = 7
$defWs.Range("E3").Value = 53
$defWs.Range("F3").Value = 181
$defWs.Range("G3").Value = 74
$defWs.Range("I3").Value = 99
$defWs.Range("J3").Value = 100
$defWs.Range("N3").Value = 32

# ---------------------------------------------------------------
# ST sheet: season totals + appended per-game logs
# ---------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 161
$stWs.Range("D2").Value = 101
$stWs.Range("F2").Value = 218
$stWs.Range("G2").Value = 205
$stWs.Range("B3").Value = 97

$stWs.Range("B4").Value2 = $stWs.Range("B4").Value2 + " 69"
$stWs.Range("B5").Value2 = $stWs.Range("B5").Value2 + " 19"
$stWs.Range("B6").Value2 = $stWs.Range("B6").Value2 + " 75"
$stWs.Range("D3").Value2 = $stWs.Range("D3").Value2 + " 41"
$stWs.Range("D4").Value2 = $stWs.Range("D4").Value2 + " 8"
$stWs.Range("D5").Value2 = $stWs.Range("D5").Value2 + " 0"

# ---------------------------------------------------------------
# TURNS sheet: season totals through this week
# ---------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B2").Value = 12
$turnsWs.Range("C2").Value = 10
$turnsWs.Range("D2").Value = 11
$turnsWs.Range("E2").Value = 11
$turnsWs.Range("E3").Value = 10

# ---------------------------------------------------------------
# PEN sheet: season totals through this week
# ---------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("D4").Value = 13
